$p = $ppt.ActivePresentation
Write-Output $p.Name
